$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.244.20'
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("D3").Value = '2.550.84'
$ws.Range("E3").Value = '  -2.39%  '

$ws.Range("E4").Value = '  +0.04%  '

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '590.66'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  +0.57%  '

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '173.59'
$cell.Style = $origStyle
$ws.Range("E6").Value = '  +4.87%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("E8").Value = '  +0.52%  '

$ws.Range("D9").Value = '2.551.64'
$ws.Range("E9").Value = '  -2.33%  '

$ws.Range("E10").Value = '  -1.15%  '

$ws.Range("E11").Value = '  +1.87%  '

$ws.Range("E12").Value = '  -0.73%  '

$ws.Range("E13").Value = '  -4.94%  '

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '26.99'
$cell.Style = $origStyle
$ws.Range("E14").Value = '  -0.98%  '

$ws.Range("D15").Value = '3.011.59'
$ws.Range("E15").Value = '  -2.53%  '

$ws.Range("E16").Value = '  -1.13%  '

$ws.Range("D17").Value = '67.107.54'
$ws.Range("E17").Value = '  +0.45%  '

$ws.Range("D18").Value = '2.559.87'
$ws.Range("E18").Value = '  -1.72%  '

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.05'
$cell.Style = $origStyle
$ws.Range("E19").Value = '  +3.30%  '

$ws.Range("E20").Value = '  -2.76%  '

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '357.10'
$cell.Style = $origStyle
$ws.Range("E21").Value = '  +0.52%  '

$ws.Range("E22").Value = '  -1.49%  '

$ws.Range("E23").Value = '  +0.67%  '

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.01'
$cell.Style = $origStyle
$ws.Range("E24").Value = '  +5.85%  '

$ws.Range("E25").Value = '  +0.02%  '

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '70.11'
$cell.Style = $origStyle
$ws.Range("E26").Value = '  +1.10%  '

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.07'
$cell.Style = $origStyle
$ws.Range("E27").Value = '  -4.21%  '

$ws.Range("E28").Value = '  -2.67%  '

$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("D30").Value = '0.0₃0990'
$ws.Range("E30").Value = '  -0.27%  '

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '533.31'
$cell.Style = $origStyle
$ws.Range("E31").Value = '  -1.40%  '

$ws.Range("E32").Value = '  +0.11%  '

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.36'
$cell.Style = $origStyle
$ws.Range("E33").Value = '  +1.33%  '

$ws.Range("E34").Value = '  -0.63%  '

$ws.Range("E35").Value = '  -0.66%  '

$ws.Range("E36").Value = '  +0.10%  '

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.47'
$cell.Style = $origStyle
$ws.Range("E37").Value = '  +0.02%  '

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '157.86'
$cell.Style = $origStyle
$ws.Range("E38").Value = '  -0.35%  '

$ws.Range("E39").Value = '  -0.88%  '

$ws.Range("E41").Value = '  -1.83%  '

$ws.Range("E42").Value = '  +0.14%  '

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.17'
$cell.Style = $origStyle
$ws.Range("E43").Value = '  +0.78%  '

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.52'
$cell.Style = $origStyle
$ws.Range("E44").Value = '  +4.80%  '

$ws.Range("E45").Value = '  +0.00%  '

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '39.72'
$cell.Style = $origStyle
$ws.Range("E46").Value = '  -1.23%  '

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '150.32'
$cell.Style = $origStyle
$ws.Range("E47").Value = '  -0.17%  '

$ws.Range("E48").Value = '  -2.18%  '

$ws.Range("D49").Value = '0.0₆0280'
$ws.Range("E49").Value = '  -4.85%  '

$ws.Range("E50").Value = '  -1.09%  '

$ws.Range("E51").Value = '  -0.06%  '

